$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3518.087
